$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap mis-ordered rows (B:AC) while keeping column A (id) fixed ---
foreach ($pair in @(@(85,86), @(117,118), @(179,180))) {
  $r1 = $pair[0]
  $r2 = $pair[1]
  $rng1 = $ws.Range("B$($r1):AC$($r1)")
  $rng2 = $ws.Range("B$($r2):AC$($r2)")
  $v1 = $rng1.Value2
  $v2 = $rng2.Value2
  $rng1.Value = $v2
  $rng2.Value = $v1
}

# --- Append new match rows 193-199 (values) ---
# Row 193
$ws.Cells.Item(193,1).Value = 191
$ws.Cells.Item(193,2).Value = 6992714
$ws.Cells.Item(193,3).Value = 'Thailand Premier League'
$ws.Cells.Item(193,4).Value = 'Thailand Premier League'
$ws.Cells.Item(193,5).Value = 45402.33333333334
$ws.Cells.Item(193,6).Value = 'BG Pathum United'
$ws.Cells.Item(193,7).Value = 'Chonburi'
$ws.Cells.Item(193,8).Value = 1
$ws.Cells.Item(193,9).Value = 1
$ws.Cells.Item(193,10).Value = 'D'
$ws.Cells.Item(193,11).Value = 1.444
$ws.Cells.Item(193,12).Value = 4.333
$ws.Cells.Item(193,13).Value = 6
$ws.Cells.Item(193,14).Value = 1.45
$ws.Cells.Item(193,15).Value = 4.5
$ws.Cells.Item(193,16).Value = 5.25
$ws.Cells.Item(193,17).Value = -1.25
$ws.Cells.Item(193,18).Value = 1.95
$ws.Cells.Item(193,19).Value = 1.85
$ws.Cells.Item(193,20).Value = 3
$ws.Cells.Item(193,21).Value = 1.95
$ws.Cells.Item(193,22).Value = 1.85
$ws.Cells.Item(193,23).Value = -1
$ws.Cells.Item(193,24).Value = 3.5
$ws.Cells.Item(193,25).Value = -1
$ws.Cells.Item(193,26).Value = -1
$ws.Cells.Item(193,27).Value = 0.8500000000000001
$ws.Cells.Item(193,28).Value = -1
$ws.Cells.Item(193,29).Value = 0.8500000000000001

# Row 194
$ws.Cells.Item(194,1).Value = 192
$ws.Cells.Item(194,2).Value = 6995900
$ws.Cells.Item(194,3).Value = 'Thailand Premier League'
$ws.Cells.Item(194,4).Value = 'Thailand Premier League'
$ws.Cells.Item(194,5).Value = 45402.35416666666
$ws.Cells.Item(194,6).Value = 'Police Tero FC'
$ws.Cells.Item(194,7).Value = 'Uthai Thani FC'
$ws.Cells.Item(194,8).Value = 1
$ws.Cells.Item(194,9).Value = 0
$ws.Cells.Item(194,10).Value = 'H'
$ws.Cells.Item(194,11).Value = 3.3
$ws.Cells.Item(194,12).Value = 3.6
$ws.Cells.Item(194,13).Value = 1.95
$ws.Cells.Item(194,14).Value = 2.7
$ws.Cells.Item(194,15).Value = 3.5
$ws.Cells.Item(194,16).Value = 2.3
$ws.Cells.Item(194,17).Value = 0.25
$ws.Cells.Item(194,18).Value = 1.775
$ws.Cells.Item(194,19).Value = 2.025
$ws.Cells.Item(194,20).Value = 3
$ws.Cells.Item(194,21).Value = 1.925
$ws.Cells.Item(194,22).Value = 1.875
$ws.Cells.Item(194,23).Value = 1.7
$ws.Cells.Item(194,24).Value = -1
$ws.Cells.Item(194,25).Value = -1
$ws.Cells.Item(194,26).Value = 0.7749999999999999
$ws.Cells.Item(194,27).Value = -1
$ws.Cells.Item(194,28).Value = -1
$ws.Cells.Item(194,29).Value = 0.875

# Row 195
$ws.Cells.Item(195,1).Value = 193
$ws.Cells.Item(195,2).Value = 6992713
$ws.Cells.Item(195,3).Value = 'Thailand Premier League'
$ws.Cells.Item(195,4).Value = 'Thailand Premier League'
$ws.Cells.Item(195,5).Value = 45402.375
$ws.Cells.Item(195,6).Value = 'Khonkaen United'
$ws.Cells.Item(195,7).Value = 'Trat FC'
$ws.Cells.Item(195,8).Value = 1
$ws.Cells.Item(195,9).Value = 1
$ws.Cells.Item(195,10).Value = 'D'
$ws.Cells.Item(195,11).Value = 2.1
$ws.Cells.Item(195,12).Value = 3.75
$ws.Cells.Item(195,13).Value = 2.875
$ws.Cells.Item(195,14).Value = 1.95
$ws.Cells.Item(195,15).Value = 3.8
$ws.Cells.Item(195,16).Value = 3.1
$ws.Cells.Item(195,17).Value = -0.25
$ws.Cells.Item(195,18).Value = 1.75
$ws.Cells.Item(195,19).Value = 1.95
$ws.Cells.Item(195,20).Value = 2.75
$ws.Cells.Item(195,21).Value = 1.825
$ws.Cells.Item(195,22).Value = 1.975
$ws.Cells.Item(195,23).Value = -1
$ws.Cells.Item(195,24).Value = 2.8
$ws.Cells.Item(195,25).Value = -1
$ws.Cells.Item(195,26).Value = -0.5
$ws.Cells.Item(195,27).Value = 0.475
$ws.Cells.Item(195,28).Value = -1
$ws.Cells.Item(195,29).Value = 0.9750000000000001

# Row 196
$ws.Cells.Item(196,1).Value = 194
$ws.Cells.Item(196,2).Value = 6992710
$ws.Cells.Item(196,3).Value = 'Thailand Premier League'
$ws.Cells.Item(196,4).Value = 'Thailand Premier League'
$ws.Cells.Item(196,5).Value = 45402.41666666666
$ws.Cells.Item(196,6).Value = 'Ratchaburi FC'
$ws.Cells.Item(196,7).Value = 'Buriram United'
$ws.Cells.Item(196,8).Value = 1
$ws.Cells.Item(196,9).Value = 4
$ws.Cells.Item(196,10).Value = 'A'
$ws.Cells.Item(196,11).Value = 5.25
$ws.Cells.Item(196,12).Value = 3.75
$ws.Cells.Item(196,13).Value = 1.571
$ws.Cells.Item(196,14).Value = 5.5
$ws.Cells.Item(196,15).Value = 3.8
$ws.Cells.Item(196,16).Value = 1.533
$ws.Cells.Item(196,17).Value = 1
$ws.Cells.Item(196,18).Value = 1.95
$ws.Cells.Item(196,19).Value = 1.85
$ws.Cells.Item(196,20).Value = 2.75
$ws.Cells.Item(196,21).Value = 1.95
$ws.Cells.Item(196,22).Value = 1.85
$ws.Cells.Item(196,23).Value = -1
$ws.Cells.Item(196,24).Value = -1
$ws.Cells.Item(196,25).Value = 0.5329999999999999
$ws.Cells.Item(196,26).Value = -1
$ws.Cells.Item(196,27).Value = 0.8500000000000001
$ws.Cells.Item(196,28).Value = 0.95
$ws.Cells.Item(196,29).Value = -1

# Row 197
$ws.Cells.Item(197,1).Value = 195
$ws.Cells.Item(197,2).Value = 6992712
$ws.Cells.Item(197,3).Value = 'Thailand Premier League'
$ws.Cells.Item(197,4).Value = 'Thailand Premier League'
$ws.Cells.Item(197,5).Value = 45403.3125
$ws.Cells.Item(197,6).Value = 'Muang Thong United'
$ws.Cells.Item(197,7).Value = 'Lamphun Warrior FC'
$ws.Cells.Item(197,8).Value = 2
$ws.Cells.Item(197,9).Value = 1
$ws.Cells.Item(197,10).Value = 'H'
$ws.Cells.Item(197,11).Value = 2
$ws.Cells.Item(197,12).Value = 3.8
$ws.Cells.Item(197,13).Value = 3
$ws.Cells.Item(197,14).Value = 1.95
$ws.Cells.Item(197,15).Value = 3.8
$ws.Cells.Item(197,16).Value = 3
$ws.Cells.Item(197,17).Value = -0.5
$ws.Cells.Item(197,18).Value = 1.975
$ws.Cells.Item(197,19).Value = 1.825
$ws.Cells.Item(197,20).Value = 2.75
$ws.Cells.Item(197,21).Value = 1.825
$ws.Cells.Item(197,22).Value = 1.975
$ws.Cells.Item(197,23).Value = 0.95
$ws.Cells.Item(197,24).Value = -1
$ws.Cells.Item(197,25).Value = -1
$ws.Cells.Item(197,26).Value = 0.9750000000000001
$ws.Cells.Item(197,27).Value = -1
$ws.Cells.Item(197,28).Value = 0.4125
$ws.Cells.Item(197,29).Value = -0.5

# Row 198
$ws.Cells.Item(198,1).Value = 196
$ws.Cells.Item(198,2).Value = 6992711
$ws.Cells.Item(198,3).Value = 'Thailand Premier League'
$ws.Cells.Item(198,4).Value = 'Thailand Premier League'
$ws.Cells.Item(198,5).Value = 45403.33333333334
$ws.Cells.Item(198,6).Value = 'Prachuap FC'
$ws.Cells.Item(198,7).Value = 'Sukhothai FC'
$ws.Cells.Item(198,8).Value = 2
$ws.Cells.Item(198,9).Value = 1
$ws.Cells.Item(198,10).Value = 'H'
$ws.Cells.Item(198,11).Value = 2.2
$ws.Cells.Item(198,12).Value = 3.5
$ws.Cells.Item(198,13).Value = 2.875
$ws.Cells.Item(198,14).Value = 1.666
$ws.Cells.Item(198,15).Value = 3.6
$ws.Cells.Item(198,16).Value = 4
$ws.Cells.Item(198,17).Value = -0.75
$ws.Cells.Item(198,18).Value = 1.925
$ws.Cells.Item(198,19).Value = 1.875
$ws.Cells.Item(198,20).Value = 2.75
$ws.Cells.Item(198,21).Value = 1.925
$ws.Cells.Item(198,22).Value = 1.875
$ws.Cells.Item(198,23).Value = 0.6659999999999999
$ws.Cells.Item(198,24).Value = -1
$ws.Cells.Item(198,25).Value = -1
$ws.Cells.Item(198,26).Value = 0.4625
$ws.Cells.Item(198,27).Value = -0.5
$ws.Cells.Item(198,28).Value = 0.4625
$ws.Cells.Item(198,29).Value = -0.5

# Row 199
$ws.Cells.Item(199,1).Value = 197
$ws.Cells.Item(199,2).Value = 6992715
$ws.Cells.Item(199,3).Value = 'Thailand Premier League'
$ws.Cells.Item(199,4).Value = 'Thailand Premier League'
$ws.Cells.Item(199,5).Value = 45403.375
$ws.Cells.Item(199,6).Value = 'Nakhon Pathom FC'
$ws.Cells.Item(199,7).Value = 'Bangkok United'
$ws.Cells.Item(199,8).Value = 1
$ws.Cells.Item(199,9).Value = 2
$ws.Cells.Item(199,10).Value = 'A'
$ws.Cells.Item(199,11).Value = 4.75
$ws.Cells.Item(199,12).Value = 3.75
$ws.Cells.Item(199,13).Value = 1.615
$ws.Cells.Item(199,14).Value = 4.75
$ws.Cells.Item(199,15).Value = 3.75
$ws.Cells.Item(199,16).Value = 1.615
$ws.Cells.Item(199,17).Value = 1
$ws.Cells.Item(199,18).Value = 1.75
$ws.Cells.Item(199,19).Value = 1.95
$ws.Cells.Item(199,20).Value = 2.75
$ws.Cells.Item(199,21).Value = 1.925
$ws.Cells.Item(199,22).Value = 1.875
$ws.Cells.Item(199,23).Value = -1
$ws.Cells.Item(199,24).Value = -1
$ws.Cells.Item(199,25).Value = 0.615
$ws.Cells.Item(199,26).Value = 0
$ws.Cells.Item(199,27).Value = 0
$ws.Cells.Item(199,28).Value = 0.4625
$ws.Cells.Item(199,29).Value = -0.5

# --- Copy cell formatting (style) for new rows from the last existing row (192) ---
$srcA = $ws.Range("A192")
$srcE = $ws.Range("E192")
$srcA.Copy() | Out-Null
$ws.Range("A193:A199").PasteSpecial(-4122) | Out-Null
$srcE.Copy() | Out-Null
$ws.Range("E193:E199").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0